# Auto-generated edit script applying scheduled market-data refresh to Jenova_Profits workbook.
# Each entry: SheetName, Row, Col, Op (set/add use Value; delete clears the cell).
$wb = $excel.ActiveWorkbook

$edits = @(
    @{Sheet="ALC"; Row=28; Col=8; Op="Set"; Value=33569.484},
    @{Sheet="ALC"; Row=28; Col=9; Op="Set"; Value=38350.223},
    @{Sheet="ALC"; Row=28; Col=10; Op="Set"; Value=1299.5},
    @{Sheet="ALC"; Row=28; Col=11; Op="Set"; Value=38350.223},
    @{Sheet="ALC"; Row=28; Col=12; Op="Set"; Value=1299.5},
    @{Sheet="ALC"; Row=28; Col=13; Op="Set"; Value=-37865.223},
    @{Sheet="ALC"; Row=28; Col=14; Op="Set"; Value=-2269.5},
    @{Sheet="ALC"; Row=43; Col=8; Op="Set"; Value=1613.1578},
    @{Sheet="ALC"; Row=43; Col=10; Op="Set"; Value=1678.5714},
    @{Sheet="ALC"; Row=43; Col=12; Op="Set"; Value=1678.5714},
    @{Sheet="ALC"; Row=43; Col=14; Op="Set"; Value=-1816.5714},
    @{Sheet="ALC"; Row=132; Col=8; Op="Set"; Value=1893.56},
    @{Sheet="ALC"; Row=132; Col=9; Op="Set"; Value=1801.7556},
    @{Sheet="ALC"; Row=132; Col=11; Op="Set"; Value=5405.266799999999},
    @{Sheet="ALC"; Row=132; Col=13; Op="Set"; Value=-2875.266799999999},
    @{Sheet="ALC"; Row=137; Col=8; Op="Set"; Value=7440.5},
    @{Sheet="ALC"; Row=137; Col=9; Op="Set"; Value=7333},
    @{Sheet="ALC"; Row=137; Col=10; Op="Set"; Value=7486.5713},
    @{Sheet="ALC"; Row=137; Col=11; Op="Set"; Value=21999},
    @{Sheet="ALC"; Row=137; Col=12; Op="Set"; Value=22459.7139},
    @{Sheet="ALC"; Row=137; Col=13; Op="Set"; Value=-19449},
    @{Sheet="ALC"; Row=137; Col=14; Op="Set"; Value=-27559.7139},
    @{Sheet="ALC"; Row=138; Col=8; Op="Set"; Value=5737.81},
    @{Sheet="ALC"; Row=138; Col=9; Op="Set"; Value=4087.8333},
    @{Sheet="ALC"; Row=138; Col=10; Op="Set"; Value=5962.8066},
    @{Sheet="ALC"; Row=138; Col=11; Op="Set"; Value=12263.4999},
    @{Sheet="ALC"; Row=138; Col=12; Op="Set"; Value=17888.4198},
    @{Sheet="ALC"; Row=138; Col=13; Op="Set"; Value=-7123.499899999999},
    @{Sheet="ALC"; Row=138; Col=14; Op="Set"; Value=-28168.4198},
    @{Sheet="ARM"; Row=32; Col=8; Op="Set"; Value=16192.429},
    @{Sheet="ARM"; Row=32; Col=9; Op="Set"; Value=8081.3115},
    @{Sheet="ARM"; Row=32; Col=10; Op="Set"; Value=29564.81},
    @{Sheet="ARM"; Row=32; Col=11; Op="Set"; Value=8081.3115},
    @{Sheet="ARM"; Row=32; Col=12; Op="Set"; Value=29564.81},
    @{Sheet="ARM"; Row=32; Col=13; Op="Set"; Value=-7794.3115},
    @{Sheet="ARM"; Row=32; Col=14; Op="Set"; Value=-30138.81},
    @{Sheet="ARM"; Row=61; Col=8; Op="Set"; Value=5314.6855},
    @{Sheet="ARM"; Row=61; Col=9; Op="Set"; Value=4625.5},
    @{Sheet="ARM"; Row=61; Col=11; Op="Set"; Value=4625.5},
    @{Sheet="ARM"; Row=61; Col=13; Op="Set"; Value=-4413.5},
    @{Sheet="ARM"; Row=125; Col=8; Op="Set"; Value=84998.5},
    @{Sheet="ARM"; Row=125; Col=9; Op="Set"; Value=40001},
    @{Sheet="ARM"; Row=125; Col=10; Op="Set"; Value=99997.664},
    @{Sheet="ARM"; Row=125; Col=11; Op="Set"; Value=40001},
    @{Sheet="ARM"; Row=125; Col=12; Op="Set"; Value=99997.664},
    @{Sheet="ARM"; Row=125; Col=13; Op="Set"; Value=-35081},
    @{Sheet="ARM"; Row=125; Col=14; Op="Set"; Value=-109837.664},
    @{Sheet="ARM"; Row=132; Col=8; Op="Set"; Value=6429.114},
    @{Sheet="ARM"; Row=132; Col=9; Op="Set"; Value=5235.25},
    @{Sheet="ARM"; Row=132; Col=11; Op="Set"; Value=15705.75},
    @{Sheet="ARM"; Row=132; Col=13; Op="Set"; Value=-13175.75},
    @{Sheet="ARM"; Row=136; Col=8; Op="Set"; Value=5314.6855},
    @{Sheet="ARM"; Row=136; Col=9; Op="Set"; Value=4625.5},
    @{Sheet="ARM"; Row=136; Col=11; Op="Set"; Value=13876.5},
    @{Sheet="ARM"; Row=136; Col=13; Op="Set"; Value=-11326.5},
    @{Sheet="BSM"; Row=11; Col=8; Op="Set"; Value=1578.375},
    @{Sheet="BSM"; Row=11; Col=9; Op="Set"; Value=85.40000000000001},
    @{Sheet="BSM"; Row=11; Col=10; Op="Set"; Value=4066.6667},
    @{Sheet="BSM"; Row=11; Col=11; Op="Set"; Value=85.40000000000001},
    @{Sheet="BSM"; Row=11; Col=12; Op="Set"; Value=4066.6667},
    @{Sheet="BSM"; Row=11; Col=13; Op="Set"; Value=54.59999999999999},
    @{Sheet="BSM"; Row=11; Col=14; Op="Set"; Value=-4346.6667},
    @{Sheet="BSM"; Row=20; Col=8; Op="Set"; Value=2042.6111},
    @{Sheet="BSM"; Row=20; Col=9; Op="Set"; Value=2072.8462},
    @{Sheet="BSM"; Row=20; Col=10; Op="Set"; Value=1964},
    @{Sheet="BSM"; Row=20; Col=11; Op="Set"; Value=2072.8462},
    @{Sheet="BSM"; Row=20; Col=12; Op="Set"; Value=1964},
    @{Sheet="BSM"; Row=20; Col=13; Op="Set"; Value=-1825.8462},
    @{Sheet="BSM"; Row=20; Col=14; Op="Set"; Value=-2458},
    @{Sheet="CRP"; Row=31; Col=8; Op="Set"; Value=1961.4445},
    @{Sheet="CRP"; Row=31; Col=9; Op="Set"; Value=1894.125},
    @{Sheet="CRP"; Row=31; Col=11; Op="Set"; Value=1894.125},
    @{Sheet="CRP"; Row=31; Col=13; Op="Set"; Value=-1599.125},
    @{Sheet="CRP"; Row=34; Col=8; Op="Set"; Value=1961.4445},
    @{Sheet="CRP"; Row=34; Col=9; Op="Set"; Value=1894.125},
    @{Sheet="CRP"; Row=34; Col=11; Op="Set"; Value=1894.125},
    @{Sheet="CRP"; Row=34; Col=13; Op="Set"; Value=-1692.125},
    @{Sheet="CRP"; Row=50; Col=8; Op="Set"; Value=22939.5},
    @{Sheet="CRP"; Row=50; Col=10; Op="Set"; Value=56000},
    @{Sheet="CRP"; Row=50; Col=12; Op="Set"; Value=56000},
    @{Sheet="CRP"; Row=50; Col=14; Op="Set"; Value=-57250},
    @{Sheet="CRP"; Row=94; Col=8; Op="Set"; Value=1346.9166},
    @{Sheet="CRP"; Row=94; Col=9; Op="Set"; Value=1333.7142},
    @{Sheet="CRP"; Row=94; Col=10; Op="Set"; Value=1365.4},
    @{Sheet="CRP"; Row=94; Col=11; Op="Set"; Value=1333.7142},
    @{Sheet="CRP"; Row=94; Col=12; Op="Set"; Value=1365.4},
    @{Sheet="CRP"; Row=94; Col=13; Op="Set"; Value=-882.7141999999999},
    @{Sheet="CRP"; Row=94; Col=14; Op="Set"; Value=-2267.4},
    @{Sheet="CRP"; Row=99; Col=8; Op="Set"; Value=5701.3335},
    @{Sheet="CRP"; Row=99; Col=9; Op="Set"; Value=5635.3335},
    @{Sheet="CRP"; Row=99; Col=10; Op="Set"; Value=5833.3335},
    @{Sheet="CRP"; Row=99; Col=11; Op="Set"; Value=5635.3335},
    @{Sheet="CRP"; Row=99; Col=12; Op="Set"; Value=5833.3335},
    @{Sheet="CRP"; Row=99; Col=13; Op="Set"; Value=-4137.3335},
    @{Sheet="CRP"; Row=99; Col=14; Op="Set"; Value=-8829.333500000001},
    @{Sheet="CRP"; Row=107; Col=8; Op="Set"; Value=649.24243},
    @{Sheet="CRP"; Row=107; Col=9; Op="Set"; Value=460.4074},
    @{Sheet="CRP"; Row=107; Col=10; Op="Set"; Value=1499},
    @{Sheet="CRP"; Row=107; Col=11; Op="Set"; Value=460.4074},
    @{Sheet="CRP"; Row=107; Col=12; Op="Set"; Value=1499},
    @{Sheet="CRP"; Row=107; Col=13; Op="Set"; Value=1459.5926},
    @{Sheet="CRP"; Row=107; Col=14; Op="Set"; Value=-5339},
    @{Sheet="CRP"; Row=126; Col=8; Op="Set"; Value=5701.3335},
    @{Sheet="CRP"; Row=126; Col=9; Op="Set"; Value=5635.3335},
    @{Sheet="CRP"; Row=126; Col=10; Op="Set"; Value=5833.3335},
    @{Sheet="CRP"; Row=126; Col=11; Op="Set"; Value=16906.0005},
    @{Sheet="CRP"; Row=126; Col=12; Op="Set"; Value=17500.0005},
    @{Sheet="CRP"; Row=126; Col=13; Op="Set"; Value=-14436.0005},
    @{Sheet="CRP"; Row=126; Col=14; Op="Set"; Value=-22440.0005},
    @{Sheet="CRP"; Row=133; Col=8; Op="Set"; Value=48532.285},
    @{Sheet="CRP"; Row=133; Col=9; Op="Set"; Value=0},
    @{Sheet="CRP"; Row=133; Col=10; Op="Set"; Value=48532.285},
    @{Sheet="CRP"; Row=133; Col=11; Op="Set"; Value=0},
    @{Sheet="CRP"; Row=133; Col=12; Op="Set"; Value=48532.285},
    @{Sheet="CRP"; Row=133; Col=13; Op="Clear"; Value=$null},
    @{Sheet="CRP"; Row=133; Col=14; Op="Set"; Value=-53592.285},
    @{Sheet="CRP"; Row=134; Col=8; Op="Set"; Value=177180.61},
    @{Sheet="CRP"; Row=134; Col=9; Op="Set"; Value=1972.1842},
    @{Sheet="CRP"; Row=134; Col=10; Op="Set"; Value=527597.5},
    @{Sheet="CRP"; Row=134; Col=11; Op="Set"; Value=5916.5526},
    @{Sheet="CRP"; Row=134; Col=12; Op="Set"; Value=1582792.5},
    @{Sheet="CRP"; Row=134; Col=13; Op="Set"; Value=-3381.5526},
    @{Sheet="CRP"; Row=134; Col=14; Op="Set"; Value=-1587862.5},
    @{Sheet="CUL"; Row=44; Col=8; Op="Set"; Value=740.1667},
    @{Sheet="CUL"; Row=44; Col=9; Op="Set"; Value=500.66666},
    @{Sheet="CUL"; Row=44; Col=10; Op="Set"; Value=979.6667},
    @{Sheet="CUL"; Row=44; Col=11; Op="Set"; Value=1501.99998},
    @{Sheet="CUL"; Row=44; Col=12; Op="Set"; Value=2939.0001},
    @{Sheet="CUL"; Row=44; Col=13; Op="Set"; Value=-1103.99998},
    @{Sheet="CUL"; Row=44; Col=14; Op="Set"; Value=-3735.0001},
    @{Sheet="CUL"; Row=107; Col=8; Op="Set"; Value=150856.42},
    @{Sheet="CUL"; Row=107; Col=9; Op="Set"; Value=1108.6666},
    @{Sheet="CUL"; Row=107; Col=10; Op="Set"; Value=263167.25},
    @{Sheet="CUL"; Row=107; Col=11; Op="Set"; Value=3325.9998},
    @{Sheet="CUL"; Row=107; Col=12; Op="Set"; Value=789501.75},
    @{Sheet="CUL"; Row=107; Col=13; Op="Set"; Value=-1405.9998},
    @{Sheet="CUL"; Row=107; Col=14; Op="Set"; Value=-793341.75},
    @{Sheet="CUL"; Row=113; Col=8; Op="Set"; Value=3088736.5},
    @{Sheet="CUL"; Row=113; Col=10; Op="Set"; Value=2527.4546},
    @{Sheet="CUL"; Row=113; Col=12; Op="Set"; Value=7582.3638},
    @{Sheet="CUL"; Row=113; Col=14; Op="Set"; Value=-11922.3638},
    @{Sheet="CUL"; Row=114; Col=8; Op="Set"; Value=126408.375},
    @{Sheet="CUL"; Row=114; Col=9; Op="Set"; Value=956},
    @{Sheet="CUL"; Row=114; Col=10; Op="Set"; Value=502765.5},
    @{Sheet="CUL"; Row=114; Col=11; Op="Set"; Value=2868},
    @{Sheet="CUL"; Row=114; Col=12; Op="Set"; Value=1508296.5},
    @{Sheet="CUL"; Row=114; Col=13; Op="Set"; Value=386},
    @{Sheet="CUL"; Row=114; Col=14; Op="Set"; Value=-1514804.5},
    @{Sheet="CUL"; Row=121; Col=8; Op="Set"; Value=1429727},
    @{Sheet="CUL"; Row=121; Col=9; Op="Set"; Value=899.5},
    @{Sheet="CUL"; Row=121; Col=10; Op="Set"; Value=2001258},
    @{Sheet="CUL"; Row=121; Col=11; Op="Set"; Value=2698.5},
    @{Sheet="CUL"; Row=121; Col=12; Op="Set"; Value=6003774},
    @{Sheet="CUL"; Row=121; Col=13; Op="Set"; Value=-1388.5},
    @{Sheet="CUL"; Row=121; Col=14; Op="Set"; Value=-6006394},
    @{Sheet="CUL"; Row=122; Col=8; Op="Set"; Value=1498.0769},
    @{Sheet="CUL"; Row=122; Col=10; Op="Set"; Value=1498.0769},
    @{Sheet="CUL"; Row=122; Col=12; Op="Set"; Value=13482.6921},
    @{Sheet="CUL"; Row=122; Col=14; Op="Set"; Value=-18382.6921},
    @{Sheet="CUL"; Row=131; Col=8; Op="Set"; Value=148719.64},
    @{Sheet="CUL"; Row=131; Col=9; Op="Set"; Value=502499.5},
    @{Sheet="CUL"; Row=131; Col=10; Op="Set"; Value=89756.336},
    @{Sheet="CUL"; Row=131; Col=11; Op="Set"; Value=1507498.5},
    @{Sheet="CUL"; Row=131; Col=12; Op="Set"; Value=269269.008},
    @{Sheet="CUL"; Row=131; Col=13; Op="Set"; Value=-1502458.5},
    @{Sheet="CUL"; Row=131; Col=14; Op="Set"; Value=-279349.008},
    @{Sheet="CUL"; Row=132; Col=8; Op="Set"; Value=1003791.94},
    @{Sheet="CUL"; Row=132; Col=10; Op="Set"; Value=1115523.1},
    @{Sheet="CUL"; Row=132; Col=12; Op="Set"; Value=10039707.9},
    @{Sheet="CUL"; Row=132; Col=14; Op="Set"; Value=-10044767.9},
    @{Sheet="CUL"; Row=134; Col=8; Op="Set"; Value=2135.348},
    @{Sheet="CUL"; Row=134; Col=9; Op="Set"; Value=1323.3182},
    @{Sheet="CUL"; Row=134; Col=11; Op="Set"; Value=3969.9546},
    @{Sheet="CUL"; Row=134; Col=13; Op="Set"; Value=1100.0454},
    @{Sheet="CUL"; Row=141; Col=8; Op="Set"; Value=6512.5},
    @{Sheet="CUL"; Row=141; Col=10; Op="Set"; Value=0},
    @{Sheet="CUL"; Row=141; Col=12; Op="Set"; Value=0},
    @{Sheet="CUL"; Row=141; Col=14; Op="Clear"; Value=$null},
    @{Sheet="GSM"; Row=80; Col=8; Op="Set"; Value=1257288.8},
    @{Sheet="GSM"; Row=80; Col=9; Op="Set"; Value=1255702.1},
    @{Sheet="GSM"; Row=80; Col=10; Op="Set"; Value=1258875.5},
    @{Sheet="GSM"; Row=80; Col=11; Op="Set"; Value=1255702.1},
    @{Sheet="GSM"; Row=80; Col=12; Op="Set"; Value=1258875.5},
    @{Sheet="GSM"; Row=80; Col=13; Op="Set"; Value=-1254704.1},
    @{Sheet="GSM"; Row=80; Col=14; Op="Set"; Value=-1260871.5},
    @{Sheet="GSM"; Row=83; Col=8; Op="Set"; Value=1257288.8},
    @{Sheet="GSM"; Row=83; Col=9; Op="Set"; Value=1255702.1},
    @{Sheet="GSM"; Row=83; Col=10; Op="Set"; Value=1258875.5},
    @{Sheet="GSM"; Row=83; Col=11; Op="Set"; Value=6278510.5},
    @{Sheet="GSM"; Row=83; Col=12; Op="Set"; Value=6294377.5},
    @{Sheet="GSM"; Row=83; Col=13; Op="Set"; Value=-6273518.5},
    @{Sheet="GSM"; Row=83; Col=14; Op="Set"; Value=-6304361.5},
    @{Sheet="GSM"; Row=117; Col=8; Op="Set"; Value=18665.334},
    @{Sheet="GSM"; Row=117; Col=10; Op="Set"; Value=18665.334},
    @{Sheet="GSM"; Row=117; Col=12; Op="Set"; Value=18665.334},
    @{Sheet="GSM"; Row=117; Col=14; Op="Set"; Value=-25549.334},
    @{Sheet="GSM"; Row=128; Col=8; Op="Set"; Value=49999.5},
    @{Sheet="GSM"; Row=128; Col=10; Op="Set"; Value=49999.5},
    @{Sheet="GSM"; Row=128; Col=12; Op="Set"; Value=49999.5},
    @{Sheet="GSM"; Row=128; Col=14; Op="Set"; Value=-59959.5},
    @{Sheet="GSM"; Row=130; Col=8; Op="Set"; Value=64995.75},
    @{Sheet="GSM"; Row=130; Col=10; Op="Set"; Value=64995.75},
    @{Sheet="GSM"; Row=130; Col=12; Op="Set"; Value=64995.75},
    @{Sheet="GSM"; Row=130; Col=14; Op="Set"; Value=-75035.75},
    @{Sheet="GSM"; Row=132; Col=8; Op="Set"; Value=30134.4},
    @{Sheet="GSM"; Row=132; Col=9; Op="Set"; Value=4956.6772},
    @{Sheet="GSM"; Row=132; Col=10; Op="Set"; Value=116857.664},
    @{Sheet="GSM"; Row=132; Col=11; Op="Set"; Value=14870.0316},
    @{Sheet="GSM"; Row=132; Col=12; Op="Set"; Value=350572.992},
    @{Sheet="GSM"; Row=132; Col=13; Op="Set"; Value=-12340.0316},
    @{Sheet="GSM"; Row=132; Col=14; Op="Set"; Value=-355632.992},
    @{Sheet="LTW"; Row=40; Col=8; Op="Set"; Value=4377.6665},
    @{Sheet="LTW"; Row=40; Col=9; Op="Set"; Value=4399.875},
    @{Sheet="LTW"; Row=40; Col=10; Op="Set"; Value=4200},
    @{Sheet="LTW"; Row=40; Col=11; Op="Set"; Value=4399.875},
    @{Sheet="LTW"; Row=40; Col=12; Op="Set"; Value=4200},
    @{Sheet="LTW"; Row=40; Col=13; Op="Set"; Value=-4263.875},
    @{Sheet="LTW"; Row=40; Col=14; Op="Set"; Value=-4472},
    @{Sheet="LTW"; Row=122; Col=8; Op="Set"; Value=8279.4},
    @{Sheet="LTW"; Row=122; Col=9; Op="Set"; Value=8199.444},
    @{Sheet="LTW"; Row=122; Col=11; Op="Set"; Value=24598.332},
    @{Sheet="LTW"; Row=122; Col=13; Op="Set"; Value=-22148.332},
    @{Sheet="WVR"; Row=100; Col=8; Op="Set"; Value=1006.5},
    @{Sheet="WVR"; Row=100; Col=9; Op="Set"; Value=1006.5},
    @{Sheet="WVR"; Row=100; Col=10; Op="Set"; Value=0},
    @{Sheet="WVR"; Row=100; Col=11; Op="Set"; Value=2013},
    @{Sheet="WVR"; Row=100; Col=12; Op="Set"; Value=0},
    @{Sheet="WVR"; Row=100; Col=13; Op="Set"; Value=-1472},
    @{Sheet="WVR"; Row=100; Col=14; Op="Clear"; Value=$null},
    @{Sheet="WVR"; Row=106; Col=8; Op="Set"; Value=79543.14},
    @{Sheet="WVR"; Row=106; Col=9; Op="Set"; Value=79491.664},
    @{Sheet="WVR"; Row=106; Col=11; Op="Set"; Value=79491.664},
    @{Sheet="WVR"; Row=106; Col=13; Op="Set"; Value=-79852},
)

foreach ($edit in $edits) {
    $ws = $wb.Worksheets.Item($edit.Sheet)
    $cell = $ws.Cells.Item($edit.Row, $edit.Col)
    if ($edit.Op -eq "Set") {
        $cell.Value2 = $edit.Value
    } else {
        $cell.ClearContents()
    }
}

Write-Host "Applied $($edits.Count) cell edits."